# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) for the affected leve rows across the eight
# crafting-job sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR), per the latest
# Universalis price pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H137").Value = 804.6
$ws.Range("I137").Value = 804.6
$ws.Range("K137").Value = 2413.8
$ws.Range("M137").Value = 136.1999999999998
$ws.Range("H138").Value = 1777.6666
$ws.Range("I138").Value = 1699.8
$ws.Range("J138").Value = 1875
$ws.Range("K138").Value = 5099.4
$ws.Range("L138").Value = 5625
$ws.Range("M138").Value = 40.60000000000036
$ws.Range("N138").Value = -15905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 799
$ws.Range("I6").Value = 599
$ws.Range("J6").Value = 999
$ws.Range("K6").Value = 599
$ws.Range("L6").Value = 999
$ws.Range("M6").Value = -426
$ws.Range("N6").Value = -1345
$ws.Range("H21").Value = 2115
$ws.Range("I21").Value = 2000
$ws.Range("J21").Value = 2345
$ws.Range("K21").Value = 2000
$ws.Range("L21").Value = 2345
$ws.Range("M21").Value = -1626
$ws.Range("N21").Value = -3093
$ws.Range("H33").Value = 25000
$ws.Range("J33").Value = 25000
$ws.Range("L33").Value = 25000
$ws.Range("N33").Value = -25658
$ws.Range("H45").Value = 4011.5
$ws.Range("I45").Value = 4011.5
$ws.Range("K45").Value = 4011.5
$ws.Range("M45").Value = -3634.5
$ws.Range("H61").Value = 881.75
$ws.Range("I61").Value = 904.3333
$ws.Range("K61").Value = 904.3333
$ws.Range("M61").Value = -692.3333
$ws.Range("H97").Value = 1168.1305
$ws.Range("I97").Value = 887.7368
$ws.Range("K97").Value = 887.7368
$ws.Range("M97").Value = -391.7368
$ws.Range("H132").Value = 13725
$ws.Range("I132").Value = 13725
$ws.Range("K132").Value = 41175
$ws.Range("M132").Value = -38645
$ws.Range("H136").Value = 881.75
$ws.Range("I136").Value = 904.3333
$ws.Range("K136").Value = 2712.9999
$ws.Range("M136").Value = -162.9998999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 5168.875
$ws.Range("I7").Value = 36.333332
$ws.Range("J7").Value = 8248.4
$ws.Range("K7").Value = 36.333332
$ws.Range("L7").Value = 8248.4
$ws.Range("M7").Value = 76.666668
$ws.Range("N7").Value = -8474.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 20001500
$ws.Range("I6").Value = 20001500
$ws.Range("K6").Value = 20001500
$ws.Range("M6").Value = -20001387
$ws.Range("H12").Value = 10000
$ws.Range("J12").Value = 10000
$ws.Range("L12").Value = 10000
$ws.Range("N12").Value = -10340
$ws.Range("H19").Value = 4286156.5
$ws.Range("I19").Value = 5000231.5
$ws.Range("K19").Value = 5000231.5
$ws.Range("M19").Value = -5000061.5
$ws.Range("H24").Value = 4286156.5
$ws.Range("I24").Value = 5000231.5
$ws.Range("K24").Value = 5000231.5
$ws.Range("M24").Value = -5000061.5
$ws.Range("H25").Value = 2183.3333
$ws.Range("J25").Value = 3000
$ws.Range("L25").Value = 3000
$ws.Range("N25").Value = -3348
$ws.Range("H32").Value = 1803.2222
$ws.Range("I32").Value = 2044.2858
$ws.Range("K32").Value = 2044.2858
$ws.Range("M32").Value = -1728.2858
$ws.Range("H99").Value = 1002199.2
$ws.Range("I99").Value = 2749
$ws.Range("K99").Value = 2749
$ws.Range("M99").Value = -1251
$ws.Range("H126").Value = 1002199.2
$ws.Range("I126").Value = 2749
$ws.Range("K126").Value = 8247
$ws.Range("M126").Value = -5777
$ws.Range("H132").Value = 2716.889
$ws.Range("I132").Value = 1806.375
$ws.Range("K132").Value = 5419.125
$ws.Range("M132").Value = -2889.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3500
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -9340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 6806
$ws.Range("I31").Value = 3010
$ws.Range("J31").Value = 12500
$ws.Range("K31").Value = 3010
$ws.Range("L31").Value = 12500
$ws.Range("M31").Value = -2762
$ws.Range("N31").Value = -12996
$ws.Range("H61").Value = 4116.6665
$ws.Range("I61").Value = 4116.6665
$ws.Range("K61").Value = 4116.6665
$ws.Range("M61").Value = -3914.6665
$ws.Range("H68").Value = 2038.8889
$ws.Range("I68").Value = 2166.6667
$ws.Range("J68").Value = 1783.3334
$ws.Range("K68").Value = 2166.6667
$ws.Range("L68").Value = 1783.3334
$ws.Range("M68").Value = -1417.6667
$ws.Range("N68").Value = -3281.3334
$ws.Range("H71").Value = 2038.8889
$ws.Range("I71").Value = 2166.6667
$ws.Range("J71").Value = 1783.3334
$ws.Range("K71").Value = 10833.3335
$ws.Range("L71").Value = 8916.666999999999
$ws.Range("M71").Value = -7089.333500000001
$ws.Range("N71").Value = -16404.667
$ws.Range("H82").Value = 1420.4166
$ws.Range("I82").Value = 1116.8334
$ws.Range("K82").Value = 1116.8334
$ws.Range("M82").Value = -755.8334
$ws.Range("H85").Value = 1420.4166
$ws.Range("I85").Value = 1116.8334
$ws.Range("K85").Value = 1116.8334
$ws.Range("M85").Value = 131.1666
$ws.Range("H113").Value = 4116.6665
$ws.Range("I113").Value = 4116.6665
$ws.Range("K113").Value = 4116.6665
$ws.Range("M113").Value = -1946.6665
$ws.Range("H122").Value = 1867
$ws.Range("I122").Value = 1867
$ws.Range("K122").Value = 5601
$ws.Range("M122").Value = -3151

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H46").Value = 99000
$ws.Range("J46").Value = 99000
$ws.Range("L46").Value = 99000
$ws.Range("N46").Value = -99462
$ws.Range("H134").Value = 99000
$ws.Range("J134").Value = 99000
$ws.Range("L134").Value = 297000
$ws.Range("N134").Value = -302070
